# Auto-committed on 2023/06/02 週五 17:31:38.49
# Add a new "findL6971" lookup entry to the DBS sheet (row 3), matching
# the pattern of the existing "findLogs" entry in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

# Write column B first, then column A, so the shared-string table gets the
# same insertion order as the authored workbook (new long text before the
# new short "findL6971" label).
$ws.Range("B3").Value = "Type = ,AND TableName = ,AND ExecuteDate = ,AND BatchNo = ,AND CustNo = ,AND FacmNo = ,AND BormNo = ,AND IsDeleted = "
$ws.Range("A3").Value = "findL6971"

# Match the row's existing formatting (style carries over automatically
# from the column style), then move the selection down to A4 as in the
# saved file.
$ws.Range("A4").Select()
